$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: the sentence about "most relevant scenario space" was previously
# split into three runs around a grammar-check proofing mark (<w:proofErr>)
# on the word "most". The edit removes that split so the whole sentence is a
# single run again. Re-running Find/Replace across the full sentence collapses
# it (and the proofErr marks) back into one run with identical visible text.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Extract all attributes and interactions related to a recommended product and identify all products linked to these attributes, to collect assumingly most relevant scenario space for a given recommended product.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Extract all attributes and interactions related to a recommended product and identify all products linked to these attributes, to collect assumingly most relevant scenario space for a given recommended product.",
    2) | Out-Null

# ---------------------------------------------------------------------------
# Change 2: "Foster" becomes "Generate options for" inside the Diversified
# Recommendations bullet. Replace just that word first …
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Foster", $true, $false, $false, $false, $false, $true, 1, $false, "Generate options for", 2) | Out-Null

# … then nudge formatting on the replacement text (set then immediately
# cleared) so the run boundaries around it are preserved, matching how Word
# keeps the newly-typed text as its own run distinct from the text before and
# after it.
$fullText = $d.Content.Text
$insIdx = $fullText.IndexOf("Generate options for")
$insLen = "Generate options for".Length
$insRange = $d.Range($insIdx, $insIdx + $insLen)
$insRange.Font.Bold = $true
$insRange.Font.Bold = $false

# ---------------------------------------------------------------------------
# Change 3: drop "not only aid users in understanding why certain products
# are recommended but also " from the Significance paragraph.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "This approach will not only aid users in understanding why certain products are recommended but also assist businesses in refining their recommendation algorithms based on user feedback and behavior analysis.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "This approach will assist businesses in refining their recommendation algorithms based on user feedback and behavior analysis.",
    2) | Out-Null
